# Website update: add the new "Honda Livo BS6" entry as the first column,
# shifting the existing Honda CB Hornet 160R / Honda XBlade / Honda CBF190R
# columns one place to the right, and refresh the price labels to the
# "Rs. <amount>" format (replacing the old "<amount> onwards" wording).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new column at A - this shifts the old A/B/C columns to B/C/D.
$ws.Columns("A").Insert()

# Populate the new first column (the Honda Livo BS6 entry).
$ws.Range("A1").Value = "Honda Livo BS6"
$ws.Range("A2").Value = "Rs. 67,000"
$ws.Range("A3").Value = "Exp. Launch : 30 Jun 2020"

# Refresh the price wording for the existing models (now in B2:D2).
$ws.Range("B2").Value = "Rs. 86,500"
$ws.Range("C2").Value = "Rs. 80,325"
$ws.Range("D2").Value = "Rs. 1.1 Lakh"
